$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some of the values below look numeric/percentage (e.g. '21', '10%') but
# must remain literal text, matching the source inlineStr cells. Force the
# NumberFormat to Text ('@') for those cells before assigning so Excel does
# not auto-convert them to numbers/percentages.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G9").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = '24EM03355'
$ws.Range("B2").Value = '21'
$ws.Range("C2").Value = 'cancer pulmonaire'
$ws.Range("D2").Value = 'CurePath'
$ws.Range("E2").Value = 'Carcinome non à petites cellules NOS'
$ws.Range("F2").Value = 'Optimale'
$ws.Range("G2").Value = '3 / 3'

# Row 3
$ws.Range("A3").Value = '24EM03456'
$ws.Range("B3").Value = '24CU052383'
$ws.Range("C3").Value = 'COSMIC'
$ws.Range("D3").Value = 'Curepath'
$ws.Range("E3").Value = 'Adénocarcinome TTF1+'
$ws.Range("F3").Value = 'Optimale'
$ws.Range("G3").Value = '<10%'

# Row 4
$ws.Range("A4").Value = '24EM03461'
$ws.Range("B4").Value = '24CU002162-4'
$ws.Range("C4").Value = 'Adénocarcinome colorectal métastatique'
$ws.Range("D4").Value = 'Adénocarcinome colorectal métastatique'
$ws.Range("E4").Value = 'Adénocarcinome colorectal métastatique'
$ws.Range("F4").Value = 'Optimale'
$ws.Range("G4").Value = '3 / 3'

# Row 5
$ws.Range("A5").Value = '24EM03462'
$ws.Range("B5").Value = '24219576'
$ws.Range("C5").Value = 'adénocarcinome mammaire'
$ws.Range("D5").Value = 'CMP'
$ws.Range("E5").Value = 'Métastase hépatique d’un adénocarcinome mammaire'
$ws.Range("F5").Value = 'Optimale'
$ws.Range("G5").Value = '3 / 3'

# Row 6
$ws.Range("A6").Value = '24EM03839'
$ws.Range("B6").Value = '24EC09559'
$ws.Range("C6").Value = '20cytologie'
$ws.Range("D6").Value = 'Erasme'
$ws.Range("E6").Value = 'PF2'
$ws.Range("F6").Value = 'Optimale'
$ws.Range("G6").Value = '3 / 3'

# Row 7
$ws.Range("A7").Value = '24EM04099'
$ws.Range("B7").Value = 'PF1'
$ws.Range("C7").Value = 'THYROID CANCER'
$ws.Range("D7").Value = 'CurePath'
$ws.Range("E7").Value = 'PF1'
$ws.Range("F7").Value = 'Optimale'
$ws.Range("G7").Value = '3 / 3'

# Row 8
$ws.Range("A8").Value = '24EM04107'
$ws.Range("B8").Value = '24CU062294-1'
$ws.Range("C8").Value = 'thyroïdiens'
$ws.Range("D8").Value = 'CurePath'
$ws.Range("E8").Value = 'PF1 oncocytaire'
$ws.Range("F8").Value = 'Optimale'
$ws.Range("G8").Value = '10%'

# Row 9
$ws.Range("A9").Value = '24EM04337'
$ws.Range("B9").Value = '8, 10'
$ws.Range("C9").Value = 'lymphomes, des cancers du sein ou d''autres cancers solides'
$ws.Range("D9").Value = 'CMP Pathology'
$ws.Range("E9").Value = 'masse gastrique'
$ws.Range("F9").Value = 'Optimale'
$ws.Range("G9").Value = '38%'

# Row 10
$ws.Range("A10").Value = '24EM04347'
$ws.Range("B10").Value = '23CU032757-1.02'
$ws.Range("C10").Value = 'carcinome urothélial invasif'
$ws.Range("D10").Value = 'carcinome urothélial invasif'
$ws.Range("E10").Value = 'carcinome urothélial invasif'
$ws.Range("F10").Value = 'Optimale'
$ws.Range("G10").Value = '18, 20'

# Row 11
$ws.Range("A11").Value = '24EM03451'
$ws.Range("B11").Value = '24BB11466 07'
$ws.Range("C11").Value = 'Tumeur de la granulosa'
$ws.Range("D11").Value = 'tumeurs de l’ovaire, de
l’endomètre et du sein'
$ws.Range("E11").Value = 'Tumeur de la granulosa'
$ws.Range("F11").Value = 'Optimale'
$ws.Range("G11").Value = '25%'

# Row 12
$ws.Range("A12").Value = '24EM03460'
$ws.Range("B12").Value = '24MH9721'
$ws.Range("C12").Value = 'COLON & LUNG CANCER'
$ws.Range("D12").Value = 'Centre Hospitalier de Mouscron'
$ws.Range("E12").Value = 'Adénocarcinome lieberkühnien'
$ws.Range("F12").Value = 'Optimale'
$ws.Range("G12").Value = '3 / 3'

# Row 13
$ws.Range("A13").Value = '24EM03308'
$ws.Range("B13").Value = '18/07/24'
$ws.Range("C13").Value = 'Adénocarcinome pulmonaire'
$ws.Range("D13").Value = 'CMP'
$ws.Range("E13").Value = 'Adénocarcinome pulmonaire'
$ws.Range("F13").Value = 'Optimale'
$ws.Range("G13").Value = '3 / 3'

# Row 14
$ws.Range("A14").Value = '24EM03352'
$ws.Range("B14").Value = '24MH9794'
$ws.Range("C14").Value = 'COLON & LUNG CANCER'
$ws.Range("D14").Value = 'Centre Hospitalier de Mouscron'
$ws.Range("E14").Value = 'Adénocarcinome lieberkühnien'
$ws.Range("F14").Value = 'Optimale'
$ws.Range("G14").Value = '3 / 3'
